# Apply updates to "上海-漫展信息.xlsx" per commit diff (gh-pages output update 456a3b4)
$wb = $excel.ActiveWorkbook

# --- Worksheet "展览" (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 2808
$ws.Range("F9").Value = 6510
$ws.Range("F11").Value = 84
$ws.Range("F12").Value = 78
$ws.Range("F13").Value = 5075
$ws.Range("F15").Value = 560
$ws.Range("F16").Value = 2666
$ws.Range("F17").Value = 1363
$ws.Range("F18").Value = 1528
$ws.Range("F19").Value = 1238
$ws.Range("F20").Value = 328
$ws.Range("F21").Value = 127
$ws.Range("D22").Value = "老沪闵路1388号舒也时代广场C栋2层 轮客行轮滑馆(闵行店)"
$ws.Range("F23").Value = 1102
$ws.Range("F24").Value = 258
$ws.Range("F25").Value = 549
$ws.Range("F26").Value = 1399
$ws.Range("F28").Value = 2120
$ws.Range("F29").Value = 596
$ws.Range("F30").Value = 45
$ws.Range("F31").Value = 43
$ws.Range("F32").Value = 116
$ws.Range("F33").Value = 264
$ws.Range("F34").Value = 1535
$ws.Range("F37").Value = 621
$ws.Range("F38").Value = 1091
$ws.Range("F41").Value = 2327
$ws.Range("F42").Value = 2590
$ws.Range("F44").Value = 151
$ws.Range("F45").Value = 5
$ws.Range("F46").Value = 281
$ws.Range("F48").Value = 114

# --- Worksheet "演出" (sheet2) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 329
$ws.Range("F12").Value = 99
$ws.Range("F17").Value = 165
$ws.Range("F18").Value = 44
$ws.Range("F27").Value = 422
$ws.Range("F36").Value = 7

# --- Worksheet "本地生活" (sheet3) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1710
$ws.Range("F8").Value = 1547
$ws.Range("F9").Value = 1824
$ws.Range("F10").Value = 2587
$ws.Range("F11").Value = 893
$ws.Range("F12").Value = 777
$ws.Range("F14").Value = 154

# --- Worksheet "全部类型" (sheet4) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F4").Value = 1710
$ws.Range("F6").Value = 2808
$ws.Range("F8").Value = 1547
$ws.Range("F10").Value = 6510
$ws.Range("F11").Value = 893
$ws.Range("F12").Value = 777
$ws.Range("F13").Value = 5075
$ws.Range("F15").Value = 560
$ws.Range("F16").Value = 2666
$ws.Range("F17").Value = 1363
$ws.Range("F18").Value = 1238
$ws.Range("F19").Value = 328
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = 127
$ws.Range("D22").Value = "老沪闵路1388号舒也时代广场C栋2层 轮客行轮滑馆(闵行店)"
$ws.Range("F23").Value = 329
$ws.Range("F24").Value = 1102
$ws.Range("F25").Value = 258
$ws.Range("F26").Value = 99
$ws.Range("F27").Value = 154
$ws.Range("F28").Value = 549
$ws.Range("F29").Value = 1399
$ws.Range("F31").Value = 2120
$ws.Range("F32").Value = 596
$ws.Range("F33").Value = 45
$ws.Range("F34").Value = 165
$ws.Range("F35").Value = 43
$ws.Range("F36").Value = 264
$ws.Range("F37").Value = 44
$ws.Range("F38").Value = 1535
$ws.Range("F40").Value = 1091
$ws.Range("F44").Value = 2327
$ws.Range("F45").Value = 2590
$ws.Range("F46").Value = 151
$ws.Range("F47").Value = 5
